$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D text cells to stay text (avoid Excel auto-numeric conversion)
# by temporarily applying a Text number format, then resetting cell style to Normal
# so no stray formatting diverges from the source file.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '59.376.02'
$c.Style = "Normal"

$ws.Range("E2").Value = '  -2.13%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.639.33'
$c.Style = "Normal"

$ws.Range("E3").Value = '  -0.25%  '

$ws.Range("E4").Value = '  +0.03%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '517.05'
$c.Style = "Normal"

$ws.Range("E5").Value = '  -1.58%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '150.25'
$c.Style = "Normal"

$ws.Range("E6").Value = '  -2.35%  '

$ws.Range("E7").Value = '  -0.31%  '

$ws.Range("E8").Value = '  -0.02%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '2.670.37'
$c.Style = "Normal"

$ws.Range("E9").Value = '  +0.44%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '6.53'
$c.Style = "Normal"

$ws.Range("E10").Value = '  +0.78%  '

$ws.Range("E11").Value = '  -1.42%  '

$ws.Range("E12").Value = '  -2.07%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.127'
$c.Style = "Normal"

$ws.Range("E13").Value = '  -0.82%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '3.101.31'
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '59.167.45'
$c.Style = "Normal"

$ws.Range("E15").Value = '  -2.52%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '21.46'
$c.Style = "Normal"

$ws.Range("E16").Value = '  -1.81%  '

$ws.Range("E17").Value = '  -1.04%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '2.661.57'
$c.Style = "Normal"

$ws.Range("E18").Value = '  +0.85%  '

$ws.Range("E19").Value = '  -2.36%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '348.11'
$c.Style = "Normal"

$ws.Range("E20").Value = '  -1.16%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '10.63'
$c.Style = "Normal"

$ws.Range("E21").Value = '  +0.32%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.26'
$c.Style = "Normal"

$ws.Range("E22").Value = '  -0.06%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"

$ws.Range("E23").Value = '  -0.30%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '61.33'
$c.Style = "Normal"

$ws.Range("E24").Value = '  +0.18%  '

$ws.Range("E25").Value = '  +0.06%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.753.14'
$c.Style = "Normal"

$ws.Range("E26").Value = '  -0.52%  '

$ws.Range("B27").Value = 'Binance-PegBSC-USD'

$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.991'
$c.Style = "Normal"

$ws.Range("E27").Value = '  -0.38%  '

$ws.Range("B28").Value = 'Kaspa'

$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.162'
$c.Style = "Normal"

$ws.Range("E28").Value = '  -3.05%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.0₃0838'
$c.Style = "Normal"

$ws.Range("E29").Value = '  -1.21%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '7.16'
$c.Style = "Normal"

$ws.Range("E30").Value = '  -1.07%  '

$ws.Range("E31").Value = '  -0.23%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '6.39'
$c.Style = "Normal"

$ws.Range("E32").Value = '  +3.90%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '19.18'
$c.Style = "Normal"

$ws.Range("E33").Value = '  -0.78%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.57'
$c.Style = "Normal"

$ws.Range("E34").Value = '  -2.62%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '149.32'
$c.Style = "Normal"

$ws.Range("E35").Value = '  -0.30%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.04'
$c.Style = "Normal"

$ws.Range("E36").Value = '  +17.42%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '4.10'
$c.Style = "Normal"

$ws.Range("E37").Value = '  +0.37%  '

$ws.Range("E38").Value = '  -1.77%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.893'
$c.Style = "Normal"

$ws.Range("E39").Value = '  -1.01%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '36.49'
$c.Style = "Normal"

$ws.Range("E40").Value = '  -0.98%  '

$ws.Range("E41").Value = '  -0.70%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '3.72'
$c.Style = "Normal"

$ws.Range("E42").Value = '  -0.80%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '294.15'
$c.Style = "Normal"

$ws.Range("E43").Value = '  -3.79%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.630'
$c.Style = "Normal"

$ws.Range("E44").Value = '  -0.93%  '

$ws.Range("E45").Value = '  -1.00%  '

$ws.Range("B46").Value = 'FirstDigitalUSD'

$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.992'
$c.Style = "Normal"

$ws.Range("E46").Value = '  -0.59%  '

$ws.Range("B47").Value = 'EnergySwap'

$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '19.81'
$c.Style = "Normal"

$ws.Range("E47").Value = '  -1.29%  '

$ws.Range("E48").Value = '  -1.80%  '

$ws.Range("E49").Value = '  -0.24%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0234'
$c.Style = "Normal"

$ws.Range("E50").Value = '  -2.35%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '19.03'
$c.Style = "Normal"

$ws.Range("E51").Value = '  -0.14%  '
